$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "Förändrad" (changed) date column C for every data row ---
$ws.Range("C2:C367").Value = 45182

# --- 2) Re-order rows 8-10 ---
# Before: row8=A 34702-2019, row9=A 10060-2022, row10=A 50909-2022
# After : row8=A 50909-2022 (updated figures), row9=A 34702-2019, row10=A 10060-2022
# Implemented as: insert a fresh row at 8 (pushes the old 8/9/10 down to 9/10/11),
# populate the new row 8 with the (updated) A 50909-2022 record, then remove the
# row that is now the duplicate A 50909-2022 record (row 11).
$ws.Rows(8).Insert()
$ws.Rows(11).Delete()

$ws.Range("A8").Value = "A 50909-2022"
$ws.Range("B8").Value = 44865
$ws.Range("C8").Value = 45182
$ws.Range("D8").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E8").Value = "SORSELE"
$ws.Range("G8").Value = 14.6
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 13
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 15
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 18
$ws.Range("R8").Value = "Grantickeporing`r`nRynkskinn`r`nGammelgransskål`r`nGarnlav`r`nGranticka`r`nGränsticka`r`nHarticka`r`nKnottrig blåslav`r`nLunglav`r`nRosenticka`r`nRödbrun blekspik`r`nTretåig hackspett`r`nUllticka`r`nVitgrynig nållav`r`nVitskaftad svartspik`r`nGulnål`r`nStuplav`r`nTrådticka"
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/artfynd/A 50909-2022.xlsx")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/kartor/A 50909-2022.png")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/klagomål/A 50909-2022.docx")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/klagomålsmail/A 50909-2022.docx")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/tillsyn/A 50909-2022.docx")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SORSELE/tillsynsmail/A 50909-2022.docx")'

# Restore the standard (non auto-fit) row height that every other data row uses.
$ws.Rows(8).RowHeight = 15

"done"
